$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) values for columns B:E
$ws.Range("B2").Value = 63.027085352698712
$ws.Range("C2").Value = 50.044056691350903
$ws.Range("D2").Value = 67.299110837388767
$ws.Range("E2").Value = 53.094941844320054

# Update row 3 (STR) values for columns B:E
$ws.Range("B3").Value = 64.919210033023205
$ws.Range("C3").Value = 45.78411335805194
$ws.Range("D3").Value = 71.873482619347598
$ws.Range("E3").Value = 45.914231469102674

# Update the active selection to reflect the edited range
$ws.Range("B1:E3").Select()
